$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows appended to the bottom of the data (dates 2021-07-26 .. 2021-08-09,
# i.e. serials 44403..44417), matching the commit "aggiornamento fino a 9 agosto 2021".
$data = @(
    @(44403, 1, 3, 76.2970498474059),
    @(44404, 0, 2, 50.8646998982706),
    @(44405, 0, 2, 50.8646998982706),
    @(44406, 0, 2, 50.8646998982706),
    @(44407, 1, 3, 76.2970498474059),
    @(44408, 0, 3, 76.2970498474059),
    @(44409, 0, 2, 50.8646998982706),
    @(44410, 0, 1, 25.4323499491353),
    @(44411, 0, 1, 25.4323499491353),
    @(44412, 0, 1, 25.4323499491353),
    @(44413, 1, 2, 50.8646998982706),
    @(44414, 2, 3, 76.2970498474059),
    @(44415, 0, 3, 76.2970498474059),
    @(44416, 1, 4, 101.7293997965412),
    @(44417, 0, 4, 101.7293997965412)
)

$startRow = 329
$endRow = $startRow + $data.Count - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Column A in the existing data uses a specific date/number style (the same
# one applied to every row above); copy it down onto the newly added rows.
$srcStyle = $ws.Range("A328")
$dstStyle = $ws.Range("A" + $startRow + ":A" + $endRow)
$srcStyle.Copy()
$dstStyle.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
